$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily regional COVID data rows (2021-08-26 through 2021-09-29),
# continuing directly after the existing last row (358, date serial 44433).
$rows = @(
    @(44434, 2303, 503, 262, 313, 1310, 754, 75, 379, 548, 656, 95, 159),
    @(44435, 2162, 383, 198, 359, 1440, 859, 86, 301, 295, 513, 136, 131),
    @(44436, 1976, 366, 187, 339, 1197, 758, 65, 256, 285, 520, 134, 92),
    @(44437, 1001, 239, 66, 215, 631, 570, 40, 204, 295, 324, 115, 141),
    @(44438, 325, 109, 44, 177, 794, 388, 28, 90, 196, 501, 16, 8),
    @(44439, 1931, 356, 160, 376, 613, 535, 58, 116, 174, 415, 85, 80),
    @(44440, 2047, 293, 109, 337, 1161, 716, 62, 202, 388, 497, 104, 104),
    @(44441, 1487, 249, 92, 335, 1414, 559, 54, 218, 374, 441, 82, 78),
    @(44442, 1224, 304, 107, 207, 1042, 525, 36, 192, 167, 411, 101, 65),
    @(44443, 1408, 232, 93, 163, 950, 469, 66, 142, 206, 427, 58, 96),
    @(44444, 676, 180, 48, 86, 1090, 361, 34, 121, 164, 173, 30, 80),
    @(44445, 265, 100, 31, 66, 930, 108, 19, 47, 92, 241, 12, 12),
    @(44446, 1541, 147, 66, 200, 634, 357, 62, 80, 193, 157, 53, 34),
    @(44447, 1241, 157, 61, 163, 1089, 428, 69, 126, 173, 339, 45, 39),
    @(44448, 1100, 146, 38, 112, 1057, 364, 22, 84, 114, 344, 20, 50),
    @(44449, 898, 113, 41, 93, 560, 366, 30, 125, 73, 287, 35, 47),
    @(44450, 1041, 174, 29, 88, 962, 312, 45, 70, 117, 214, 30, 31),
    @(44451, 535, 110, 17, 48, 637, 219, 18, 49, 72, 162, 47, 39),
    @(44452, 177, 26, 9, 38, 371, 49, 11, 40, 78, 104, 6, 8),
    @(44453, 1219, 107, 49, 114, 635, 261, 38, 48, 89, 144, 47, 34),
    @(44454, 1008, 131, 34, 95, 637, 286, 32, 54, 110, 198, 44, 13),
    @(44455, 819, 100, 30, 93, 666, 280, 48, 46, 75, 175, 59, 41),
    @(44456, 746, 58, 29, 92, 844, 268, 23, 39, 83, 168, 39, 23),
    @(44457, 725, 90, 23, 89, 810, 224, 35, 53, 62, 135, 29, 16),
    @(44458, 296, 40, 9, 44, 761, 143, 17, 55, 83, 66, 28, 13),
    @(44459, 89, 52, 7, 28, 321, 27, 16, 23, 46, 78, 2, 4),
    @(44460, 810, 75, 30, 65, 419, 159, 30, 17, 44, 149, 27, 23),
    @(44461, 653, 75, 17, 65, 416, 190, 23, 31, 69, 86, 58, 19),
    @(44462, 569, 89, 23, 35, 430, 133, 18, 33, 49, 166, 11, 27),
    @(44463, 541, 50, 10, 39, 565, 118, 22, 28, 48, 149, 32, 18),
    @(44464, 549, 58, 17, 51, 433, 113, 20, 24, 31, 106, 14, 28),
    @(44465, 235, 35, 13, 17, 255, 93, 12, 21, 13, 34, 3, 3),
    @(44466, 80, 9, 6, 7, 201, 19, 20, 8, 11, 31, 2, 0),
    @(44467, 527, 33, 12, 39, 302, 98, 22, 6, 37, 90, 16, 10),
    @(44468, 447, 32, 6, 29, 298, 120, 21, 31, 34, 38, 12, 14),
)

$startRow = 359
$lastExistingRow = 358
$endRow = $startRow + $rows.Count - 1

# Copy the style/format (border, bold date column, number format) of the
# last existing row down onto the new block before writing values, so the
# new rows visually match the rest of the table (date column keeps its
# date number format / border style).
$ws.Range("A" + $lastExistingRow + ":M" + $lastExistingRow).Copy()
$ws.Range("A" + $startRow + ":M" + $endRow).PasteSpecial(-4122)

$r = $startRow
foreach ($row in $rows) {
    for ($c = 1; $c -le $row.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
    $r++
}

# Update the view so it matches where the author left the selection after
# entering the new data: scrolled down near the bottom of the sheet, with
# M386 as the active cell.
$ws.Range("M386").Select()
$excel.ActiveWindow.ScrollRow = 373
$excel.ActiveWindow.ScrollColumn = 1
